$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "josef@gmail"
$ws.Range("B3").Value = "B@positive5"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:josef@gmail")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:B@positive5")
$ws.Range("D13").Select()
